$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1028208879615562
$ws.Range("D2").Value = 0.02294331548571904
$ws.Range("E2").Value = 0.1597174178937735
$ws.Range("F2").Value = 0.4349732476443577
$ws.Range("G2").Value = 0.2809637316095177
$ws.Range("H2").Value = 0.4531015744136511
$ws.Range("I2").Value = 0.4206328855603294
$ws.Range("K2").Value = 0.6004826738355575
$ws.Range("M2").Value = 0.2605831874329567
$ws.Range("O2").Value = 1.391853284685638
$ws.Range("B3").Value = 0.0908534429512855
$ws.Range("D3").Value = 0.0202945348015291
$ws.Range("E3").Value = 0.1526160146998379
$ws.Range("F3").Value = 0.4334557073628815
$ws.Range("G3").Value = 0.2807243280695673
$ws.Range("H3").Value = 0.4565789266322682
$ws.Range("I3").Value = 0.4280947762624741
$ws.Range("K3").Value = 0.5290296786291719
$ws.Range("M3").Value = 0.2301024919045318
$ws.Range("O3").Value = 1.398379298319711
$ws.Range("B4").Value = 0.08349532691781292
$ws.Range("D4").Value = 0.01865918263953859
$ws.Range("E4").Value = 0.1484061344740084
$ws.Range("F4").Value = 0.4328665598086232
$ws.Range("G4").Value = 0.2808488537790623
$ws.Range("H4").Value = 0.4589598735723328
$ws.Range("I4").Value = 0.4329459659046506
$ws.Range("K4").Value = 0.4849356880245637
$ws.Range("M4").Value = 0.2113997324883883
$ws.Range("O4").Value = 1.403467984380228
$ws.Range("B5").Value = 0.08049454096725128
$ws.Range("D5").Value = 0.01799054525110932
$ws.Range("E5").Value = 0.1467281407099392
$ws.Range("F5").Value = 0.4327125578546571
$ws.Range("G5").Value = 0.2809677101476566
$ws.Range("H5").Value = 0.4599919470469089
$ws.Range("I5").Value = 0.4349906256195148
$ws.Range("K5").Value = 0.4669125276632258
$ws.Range("M5").Value = 0.2037814901906216
$ws.Range("O5").Value = 1.405813309202372
$ws.Range("B6").Value = 0.07999613113018711
$ws.Range("D6").Value = 0.01787938580648785
$ws.Range("E6").Value = 0.1464517725870564
$ws.Range("F6").Value = 0.432692182754586
$ws.Range("G6").Value = 0.2809915550580158
$ws.Range("H6").Value = 0.4601670557400084
$ws.Range("I6").Value = 0.4353342317340179
$ws.Range("K6").Value = 0.4639165396782232
$ws.Range("M6").Value = 0.2025166893212358
$ws.Range("O6").Value = 1.406219144804183
$ws.Range("B7").Value = 0.08345486624457976
$ws.Range("D7").Value = 0.01865017408599101
$ws.Range("E7").Value = 0.1483833526949923
$ws.Range("F7").Value = 0.4328641344497015
$ws.Range("G7").Value = 0.2808501811579589
$ws.Range("H7").Value = 0.4589735421763024
$ws.Range("I7").Value = 0.4329732665695594
$ws.Range("K7").Value = 0.4846928404423068
$ws.Range("M7").Value = 0.2112969767269206
$ws.Range("O7").Value = 1.403498514848252
$ws.Range("B8").Value = 0.09869673702280579
$ws.Range("D8").Value = 0.02203190653815312
$ws.Range("E8").Value = 0.1572374804969456
$ws.Range("F8").Value = 0.4343788396522257
$ws.Range("G8").Value = 0.2808247311294139
$ws.Range("H8").Value = 0.4542495415058596
$ws.Range("I8").Value = 0.4231497860354567
$ws.Range("K8").Value = 0.5758923711019577
$ws.Range("M8").Value = 0.2500709203067615
$ws.Range("O8").Value = 1.393878733997113
$ws.Range("B9").Value = 0.1284973507118679
$ws.Range("D9").Value = 0.0285905711817378
$ws.Range("E9").Value = 0.1758063890879882
$ws.Range("F9").Value = 0.4400721224626807
$ws.Range("G9").Value = 0.2829371628452151
$ws.Range("H9").Value = 0.446936576954208
$ws.Range("I9").Value = 0.4060255309874505
$ws.Range("K9").Value = 0.7529314564872038
$ws.Range("M9").Value = 0.3262035431604247
$ws.Range("O9").Value = 1.383614539484284
$ws.Range("B10").Value = 0.1503280676075747
$ws.Range("D10").Value = 0.03336310864197856
$ws.Range("E10").Value = 0.1902030479960715
$ws.Range("F10").Value = 0.4459222927189046
$ws.Range("G10").Value = 0.2858191153492129
$ws.Range("H10").Value = 0.4427534781430751
$ws.Range("I10").Value = 0.3947493094137942
$ws.Range("K10").Value = 0.8818548959150405
$ws.Range("M10").Value = 0.3822000098169411
$ws.Range("O10").Value = 1.381342580843352
$ws.Range("B11").Value = 0.1602436500909192
$ws.Range("D11").Value = 0.03552391845624925
$ws.Range("E11").Value = 0.1969203447432335
$ws.Range("F11").Value = 0.4489473689751975
$ws.Range("G11").Value = 0.2874215255791626
$ws.Range("H11").Value = 0.4411089459119353
$ws.Range("I11").Value = 0.3899029553995907
$ws.Range("K11").Value = 0.9402470342642459
$ws.Range("M11").Value = 0.4076886010403484
$ws.Range("O11").Value = 1.38145890488326
$ws.Range("B12").Value = 0.163996012638961
$ws.Range("D12").Value = 0.0363406504953474
$ws.Range("E12").Value = 0.1994884954958493
$ws.Range("F12").Value = 0.4501453049877924
$ws.Range("G12").Value = 0.28807040806619
$ws.Range("H12").Value = 0.4405233680263052
$ws.Range("I12").Value = 0.3881085240978766
$ws.Range("K12").Value = 0.9623207818058859
$ws.Range("M12").Value = 0.4173426438455436
$ws.Range("O12").Value = 1.381668726217185
$ws.Range("B13").Value = 0.1631879873116162
$ws.Range("D13").Value = 0.03616482097658036
$ws.Range("E13").Value = 0.1989343062618047
$ws.Range("F13").Value = 0.4498849758677608
$ws.Range("G13").Value = 0.2879287850065424
$ws.Range("H13").Value = 0.4406478292673484
$ws.Range("I13").Value = 0.3884931726479386
$ws.Range("K13").Value = 0.9575685200358635
$ws.Range("M13").Value = 0.41526338370997
$ws.Range("O13").Value = 1.381616157953545
$ws.Range("B14").Value = 0.1605524098143007
$ws.Range("D14").Value = 0.03559114221984316
$ws.Range("E14").Value = 0.1971311363570791
$ws.Range("F14").Value = 0.4490448730982735
$ws.Range("G14").Value = 0.2874740650570828
$ws.Range("H14").Value = 0.4410600250173644
$ws.Range("I14").Value = 0.3897545089446761
$ws.Range("K14").Value = 0.9420638255623714
$ws.Range("M14").Value = 0.4084828044802578
$ws.Range("O14").Value = 1.381472842271108
$ws.Range("B15").Value = 0.1589377155893885
$ws.Range("D15").Value = 0.03523954810584939
$ws.Range("E15").Value = 0.1960298357060211
$ws.Range("F15").Value = 0.4485371132732396
$ws.Range("G15").Value = 0.2872010219862347
$ws.Range("H15").Value = 0.4413173480225794
$ws.Range("I15").Value = 0.3905324256865623
$ws.Range("K15").Value = 0.9325617536523509
$ws.Range("M15").Value = 0.4043297681443647
$ws.Range("O15").Value = 1.38140665823633
$ws.Range("B16").Value = 0.1496797312675255
$ws.Range("D16").Value = 0.03322168453372853
$ws.Range("E16").Value = 0.1897674643758265
$ws.Range("F16").Value = 0.445731925331927
$ws.Range("G16").Value = 0.2857202733019335
$ws.Range("H16").Value = 0.442866152835407
$ws.Range("I16").Value = 0.3950717337977614
$ws.Range("K16").Value = 0.87803357350532
$ws.Range("M16").Value = 0.3805345685856594
$ws.Range("O16").Value = 1.381358153107129
$ws.Range("B17").Value = 0.1439961554062563
$ws.Range("D17").Value = 0.03198113403332314
$ws.Range("E17").Value = 0.1859689717711532
$ws.Range("F17").Value = 0.4441042774382922
$ws.Range("G17").Value = 0.2848866454176289
$ws.Range("H17").Value = 0.4438824832451758
$ws.Range("I17").Value = 0.3979290334934285
$ws.Range("K17").Value = 0.8445158574663481
$ws.Range("M17").Value = 0.3659408413753056
$ws.Range("O17").Value = 1.381623204216339
$ws.Range("B18").Value = 0.140725685779671
$ws.Range("D18").Value = 0.03126664103413646
$ws.Range("E18").Value = 0.183799993393535
$ws.Range("F18").Value = 0.4432023350968848
$ws.Range("G18").Value = 0.2844345814359031
$ws.Range("H18").Value = 0.4444913686483645
$ws.Range("I18").Value = 0.3995991353114978
$ws.Range("K18").Value = 0.8252133561054507
$ws.Range("M18").Value = 0.3575483859152371
$ws.Range("O18").Value = 1.381883858662576
$ws.Range("B19").Value = 0.1396181234974421
$ws.Range("D19").Value = 0.03102456233528983
$ws.Range("E19").Value = 0.1830683223905822
$ws.Range("F19").Value = 0.4429028304847122
$ws.Range("G19").Value = 0.2842862236155526
$ws.Range("H19").Value = 0.4447017031143048
$ws.Range("I19").Value = 0.4001691803368033
$ws.Range("K19").Value = 0.8186737791294547
$ws.Range("M19").Value = 0.3547071022779917
$ws.Range("O19").Value = 1.381990682932212
$ws.Range("B20").Value = 0.144601331277542
$ws.Range("D20").Value = 0.0321132925695693
$ws.Range("E20").Value = 0.1863716891261689
$ws.Range("F20").Value = 0.444273999202224
$ws.Range("G20").Value = 0.2849725476759914
$ws.Range("H20").Value = 0.4437717760926745
$ws.Range("I20").Value = 0.397622109665428
$ws.Range("K20").Value = 0.8480863687086924
$ws.Range("M20").Value = 0.3674942175625731
$ws.Range("O20").Value = 1.381583787268738
$ws.Range("B21").Value = 0.1613266113845384
$ws.Range("D21").Value = 0.03575968719670186
$ws.Range("E21").Value = 0.1976601051691489
$ws.Range("F21").Value = 0.4492902086670085
$ws.Range("G21").Value = 0.2876064837224561
$ws.Range("H21").Value = 0.4409379442022612
$ws.Range("I21").Value = 0.3893829168341871
$ws.Range("K21").Value = 0.946618974248338
$ws.Range("M21").Value = 0.4104743701930147
$ws.Range("O21").Value = 1.381510435186698
$ws.Range("B22").Value = 0.1722431335039687
$ws.Range("D22").Value = 0.03813393686257882
$ws.Range("E22").Value = 0.2051804132196295
$ws.Range("F22").Value = 0.4528740905994368
$ws.Range("G22").Value = 0.2895732974638605
$ws.Range("H22").Value = 0.4393025383601952
$ws.Range("I22").Value = 0.3842358336755511
$ws.Range("K22").Value = 1.010793120207552
$ws.Range("M22").Value = 0.4385763532617091
$ws.Range("O22").Value = 1.382428936076565
$ws.Range("B23").Value = 0.1664181816025234
$ws.Range("D23").Value = 0.03686758318411876
$ws.Range("E23").Value = 0.201153538805336
$ws.Range("F23").Value = 0.4509333226352226
$ws.Range("G23").Value = 0.2885010594498283
$ws.Range("H23").Value = 0.4401555549180358
$ws.Range("I23").Value = 0.3869611634883561
$ws.Range("K23").Value = 0.9765629695094162
$ws.Range("M23").Value = 0.4235767552937517
$ws.Range("O23").Value = 1.381850145884698
$ws.Range("B24").Value = 0.1443277402998433
$ws.Range("D24").Value = 0.03205354769110613
$ws.Range("E24").Value = 0.1861895744519444
$ws.Range("F24").Value = 0.4441971626712089
$ws.Range("G24").Value = 0.284933626553169
$ws.Range("H24").Value = 0.443821750220863
$ws.Range("I24").Value = 0.3977607845762203
$ws.Range("K24").Value = 0.8464722423972546
$ws.Range("M24").Value = 0.3667919434285665
$ws.Range("O24").Value = 1.381601270438296
$ws.Range("B25").Value = 0.1204459984245432
$ws.Range("D25").Value = 0.02682425620974271
$ws.Range("E25").Value = 0.1706519454576849
$ws.Range("F25").Value = 0.4382396649420883
$ws.Range("G25").Value = 0.282132983732609
$ws.Range("H25").Value = 0.448706032422038
$ws.Range("I25").Value = 0.4104289768271001
$ws.Range("K25").Value = 0.7052357031759016
$ws.Range("M25").Value = 0.3055969885380421
$ws.Range("O25").Value = 1.385467845242147
